# Fix the last remaining "chinook" -> "steelhead" mislabeling in rows 37-41
# (carter_2005_temp_thresholds.xlsx), and update the saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($r in 37..41) {
    $ws.Cells.Item($r, 1).Value = "steelhead"
}

# Update the view so the saved selection/scroll matches the target file.
$ws.Activate()
$ws.Range("B40").Select()
$excel.ActiveWindow.ScrollRow = 28
